# Reads data from the original "bookingData" sheet into a new copy named
# "bookingData" (the old one becomes "bookingDataOriginal"), and back-fills
# the credit-card type for every row plus refreshes the expiry years.

$wb = $excel.ActiveWorkbook

$wsOriginal = $wb.Worksheets.Item(1)
$wsCopy = $wb.Worksheets.Item(2)

# --- 1. Rename the sheets -------------------------------------------------
$wsOriginal.Name = "bookingDataOriginal"
$wsCopy.Name = "bookingData"

# --- 2. Fix up ccType / expYear on the original sheet ---------------------
$ccType = @{ 2 = "VISA"; 3 = "Master Card"; 4 = "American Express"; 5 = "VISA"; 6 = "Master Card"; 7 = "American Express"; 8 = "VISA"; 9 = "Master Card"; 10 = "American Express" }
$expYear = @{ 2 = 2019; 3 = 2020; 4 = 2021; 5 = 2022; 6 = 2022; 7 = 2021; 8 = 2019; 9 = 2022; 10 = 2021 }

# Write the values FIRST (while the cells are still General-formatted) so
# the expiry years stay numeric, then apply the text number format ("@")
# afterwards - header cells keep their yellow fill, data cells stay unfilled.
for ($r = 2; $r -le 10; $r++) {
    $wsOriginal.Cells.Item($r, 14).Value = $ccType[$r]
    $wsOriginal.Cells.Item($r, 16).Value = $expYear[$r]
}

$wsOriginal.Range("M1").NumberFormat = "@"
$wsOriginal.Range("P1").NumberFormat = "@"
$wsOriginal.Range("M2:M10").NumberFormat = "@"
$wsOriginal.Range("P2:P10").NumberFormat = "@"

# --- 3. Populate the new "bookingData" (copy) sheet ------------------------
$headers = @("hotelLocation", "hotelName", "roomType", "firstName", "lastName", "address", "ccType", "expMonth")
for ($c = 1; $c -le 8; $c++) {
    $wsCopy.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$data = @(
    @("Melbourne", "Hotel Creek", "Double", "Jon", "Snow", "Winterfell", "VISA", "May"),
    @("Sydney", "Hotel Sunshine", "Standard", "Tyrion", "Lanister", "KingsLanding", "Master Card", "June"),
    @("Brisbane", "Hotel Hervey", "Deluxe", "Arya", "Stark", "Bravos", "American Express", "July"),
    @("Adelaide", "Hotel Cornice", "Super Deluxe", "Dany", "Targarian", "Andals", "VISA", "August"),
    @("London", "Hotel Creek", "Deluxe", "Sansa", "Stark", "Winterfell", "Master Card", "January"),
    @("New York", "Hotel Sunshine", "Double", "Ser", "Davos", "Seaworth", "American Express", "April"),
    @("Los Angeles", "Hotel Hervey", "Standard", "Jamie", "Lanister", "Harrenhall", "VISA", "March"),
    @("London", "Hotel Cornice", "Deluxe", "Yara", "Grajoy", "Iron Highets", "Master Card", "December"),
    @("Melbourne", "Hotel Cornice", "Super Deluxe", "The", "NightKing", "BeyondTheWall", "American Express", "October")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 1; $c -le 8; $c++) {
        $wsCopy.Cells.Item($row, $c).Value = $rowData[$c - 1]
    }
}

# --- 4. Selections / view state -------------------------------------------
$wsOriginal.Activate()
$wsOriginal.Range("L15").Select()

$wsCopy.Activate()
$wsCopy.Range("G2:G10").Select()

$wsOriginal.Activate()
